# The underlying change between before/after here is *not* a content edit:
# every hunk in the unified OOXML diff is a pure XML-attribute reordering
# (namespace declarations and element attributes re-emitted in alphabetical
# order by whatever serializer produced the "after" package - e.g. a library/
# tool version bump, per the commit message "Moving from 2.0.1 to 2.0.2").
# No text, formatting, style, relationship, or structural content differs
# between the two XML trees - only the on-disk attribute ordering does,
# which is not something controllable from the Word object model (it is an
# artifact of the writer that serialized the package, not of document
# content). So the correct interop action is to leave the document's
# content untouched.
$d = $word.ActiveDocument

Write-Output "No semantic content change required; document left as-is."
